$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67634766666667
$ws.Range("H2").Value = 83.029043
$ws.Range("I2").Value = 0.005965811625935536
$ws.Range("J2").Value = 0.005965811625935536
$ws.Range("M2").Value = 0.4216986666666666
$ws.Range("N2").Value = 1.265096
$ws.Range("O2").Value = 0.2516921781598699
$ws.Range("P2").Value = 0.2516921781598699
$ws.Range("Q2").Value = 11.67107890923644
$ws.Range("R2").Value = 105.039710183128
$ws.Range("S2").Value = 0.00150154812262319
$ws.Range("T2").Value = 0.00150154812262319
$ws.Range("G3").Value = 27.67634766666667
$ws.Range("H3").Value = 83.029043
$ws.Range("I3").Value = 0.005965811625935536
$ws.Range("J3").Value = 0.005965811625935536
$ws.Range("O3").Value = 0.2613966125002536
$ws.Range("P3").Value = 0.2613966125002536
$ws.Range("Q3").Value = 12.121077871398
$ws.Range("R3").Value = 109.089700842582
$ws.Range("S3").Value = 0.001559442949834179
$ws.Range("T3").Value = 0.001559442949834179
$ws.Range("G4").Value = 27.67634766666667
$ws.Range("H4").Value = 83.029043
$ws.Range("I4").Value = 0.005965811625935536
$ws.Range("J4").Value = 0.005965811625935536
$ws.Range("M4").Value = 0.4328273333333333
$ws.Range("N4").Value = 1.298482
$ws.Range("O4").Value = 0.2583343579312433
$ws.Range("P4").Value = 0.2583343579312433
$ws.Range("Q4").Value = 11.97907975696955
$ws.Range("R4").Value = 107.811717812726
$ws.Range("S4").Value = 0.001541174115924803
$ws.Range("T4").Value = 0.001541174115924803
$ws.Range("G5").Value = 27.67634766666667
$ws.Range("H5").Value = 83.029043
$ws.Range("I5").Value = 0.005965811625935536
$ws.Range("J5").Value = 0.005965811625935536
$ws.Range("M5").Value = 0.38297
$ws.Range("N5").Value = 1.14891
$ws.Range("O5").Value = 0.2285768514086331
$ws.Range("P5").Value = 0.2285768514086331
$ws.Range("Q5").Value = 10.59921086590333
$ws.Range("R5").Value = 95.39289779312999
$ws.Range("S5").Value = 0.001363646437553363
$ws.Range("T5").Value = 0.001363646437553363
$ws.Range("I6").Value = 0.009118181457976757
$ws.Range("J6").Value = 0.009118181457976757
$ws.Range("M6").Value = 0.4216986666666666
$ws.Range("N6").Value = 1.265096
$ws.Range("O6").Value = 0.2516921781598699
$ws.Range("P6").Value = 0.2516921781598699
$ws.Range("Q6").Value = 17.83814541547733
$ws.Range("R6").Value = 160.543308739296
$ws.Range("S6").Value = 0.002294974952015108
$ws.Range("T6").Value = 0.002294974952015108
$ws.Range("I7").Value = 0.009118181457976757
$ws.Range("J7").Value = 0.009118181457976757
$ws.Range("O7").Value = 0.2613966125002536
$ws.Range("P7").Value = 0.2613966125002536
$ws.Range("S7").Value = 0.002383461745277748
$ws.Range("T7").Value = 0.002383461745277748
$ws.Range("I8").Value = 0.009118181457976757
$ws.Range("J8").Value = 0.009118181457976757
$ws.Range("M8").Value = 0.4328273333333333
$ws.Range("N8").Value = 1.298482
$ws.Range("O8").Value = 0.2583343579312433
$ws.Range("P8").Value = 0.2583343579312433
$ws.Range("Q8").Value = 18.30889571651466
$ws.Range("R8").Value = 164.780061448632
$ws.Range("S8").Value = 0.002355539552446993
$ws.Range("T8").Value = 0.002355539552446993
$ws.Range("I9").Value = 0.009118181457976757
$ws.Range("J9").Value = 0.009118181457976757
$ws.Range("M9").Value = 0.38297
$ws.Range("N9").Value = 1.14891
$ws.Range("O9").Value = 0.2285768514086331
$ws.Range("P9").Value = 0.2285768514086331
$ws.Range("Q9").Value = 16.19989601524
$ws.Range("R9").Value = 145.79906413716
$ws.Range("S9").Value = 0.002084205208236906
$ws.Range("T9").Value = 0.002084205208236906
$ws.Range("G10").Value = 29.593002
$ws.Range("H10").Value = 88.779006
$ws.Range("I10").Value = 0.006378958578792732
$ws.Range("J10").Value = 0.006378958578792732
$ws.Range("M10").Value = 0.4216986666666666
$ws.Range("N10").Value = 1.265096
$ws.Range("O10").Value = 0.2516921781598699
$ws.Range("P10").Value = 0.2516921781598699
$ws.Range("Q10").Value = 12.479329486064
$ws.Range("R10").Value = 112.313965374576
$ws.Range("S10").Value = 0.001605533979087931
$ws.Range("T10").Value = 0.001605533979087931
$ws.Range("G11").Value = 29.593002
$ws.Range("H11").Value = 88.779006
$ws.Range("I11").Value = 0.006378958578792732
$ws.Range("J11").Value = 0.006378958578792732
$ws.Range("O11").Value = 0.2613966125002536
$ws.Range("P11").Value = 0.2613966125002536
$ws.Range("Q11").Value = 12.960491969916
$ws.Range("R11").Value = 116.644427729244
$ws.Range("S11").Value = 0.001667438163775853
$ws.Range("T11").Value = 0.001667438163775853
$ws.Range("G12").Value = 29.593002
$ws.Range("H12").Value = 88.779006
$ws.Range("I12").Value = 0.006378958578792732
$ws.Range("J12").Value = 0.006378958578792732
$ws.Range("M12").Value = 0.4328273333333333
$ws.Range("N12").Value = 1.298482
$ws.Range("O12").Value = 0.2583343579312433
$ws.Range("P12").Value = 0.2583343579312433
$ws.Range("Q12").Value = 12.808660140988
$ws.Range("R12").Value = 115.277941268892
$ws.Range("S12").Value = 0.001647904168722416
$ws.Range("T12").Value = 0.001647904168722417
$ws.Range("G13").Value = 29.593002
$ws.Range("H13").Value = 88.779006
$ws.Range("I13").Value = 0.006378958578792732
$ws.Range("J13").Value = 0.006378958578792732
$ws.Range("M13").Value = 0.38297
$ws.Range("N13").Value = 1.14891
$ws.Range("O13").Value = 0.2285768514086331
$ws.Range("P13").Value = 0.2285768514086331
$ws.Range("Q13").Value = 11.33323197594
$ws.Range("R13").Value = 101.99908778346
$ws.Range("S13").Value = 0.001458082267206532
$ws.Range("T13").Value = 0.001458082267206532
$ws.Range("G14").Value = 4539.588785666667
$ws.Range("H14").Value = 13618.766357
$ws.Range("I14").Value = 0.9785370483372949
$ws.Range("J14").Value = 0.978537048337295
$ws.Range("M14").Value = 0.4216986666666666
$ws.Range("N14").Value = 1.265096
$ws.Range("O14").Value = 0.2516921781598699
$ws.Range("P14").Value = 0.2516921781598699
$ws.Range("Q14").Value = 1914.338538130586
$ws.Range("R14").Value = 17229.04684317527
$ws.Range("S14").Value = 0.2462901211061437
$ws.Range("T14").Value = 0.2462901211061437
$ws.Range("G15").Value = 4539.588785666667
$ws.Range("H15").Value = 13618.766357
$ws.Range("I15").Value = 0.9785370483372949
$ws.Range("J15").Value = 0.978537048337295
$ws.Range("O15").Value = 0.2613966125002536
$ws.Range("P15").Value = 0.2613966125002536
$ws.Range("Q15").Value = 1988.149225393002
$ws.Range("R15").Value = 17893.34302853702
$ws.Range("S15").Value = 0.2557862696413659
$ws.Range("T15").Value = 0.2557862696413659
$ws.Range("G16").Value = 4539.588785666667
$ws.Range("H16").Value = 13618.766357
$ws.Range("I16").Value = 0.9785370483372949
$ws.Range("J16").Value = 0.978537048337295
$ws.Range("M16").Value = 0.4328273333333333
$ws.Range("N16").Value = 1.298482
$ws.Range("O16").Value = 0.2583343579312433
$ws.Range("P16").Value = 0.2583343579312433
$ws.Range("Q16").Value = 1964.858108530008
$ws.Range("R16").Value = 17683.72297677007
$ws.Range("S16").Value = 0.2527897400941491
$ws.Range("T16").Value = 0.2527897400941491
$ws.Range("G17").Value = 4539.588785666667
$ws.Range("H17").Value = 13618.766357
$ws.Range("I17").Value = 0.9785370483372949
$ws.Range("J17").Value = 0.978537048337295
$ws.Range("M17").Value = 0.38297
$ws.Range("N17").Value = 1.14891
$ws.Range("O17").Value = 0.2285768514086331
$ws.Range("P17").Value = 0.2285768514086331
$ws.Range("Q17").Value = 1738.526317246763
$ws.Range("R17").Value = 15646.73685522087
$ws.Range("S17").Value = 0.2236709174956363
$ws.Range("T17").Value = 0.2236709174956363
